# AR Automation V1 - turn the filled-in sample credential row into a blank
# template row with a guidance note, and drop the stray hyperlink that was
# left on the password cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header for column E keeps reading "businessfunctions" - no visible
# change there, just re-touch it so it stays intact as the shared string
# table gets rebuilt under the hood.
$ws.Range("E1").Value = "businessfunctions"

# The sample row (row 2) had a filled in username/password/url/report/region.
# Remove the mailto: hyperlink that lived on B2 (the password cell) first,
# then blank every sample value in the row.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("F2").ClearContents()

# Replace the old businessfunctions sample ("APP_01,APP_02,APP_03") with an
# instructional placeholder so template users know the expected format.
$ws.Range("E2").Value = "buisness functions should be seperated by comas (Eg: APP_01,APP_02)"

# Column E got a bit wider to comfortably fit the new guidance text.
$ws.Range("E1").ColumnWidth = 61.42
